# Re-order the stock rows (A2:B10) to reflect the new aggregation order.
# The Variant -> Stock mapping itself is unchanged; only the row order differs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B must keep storing its values as text (e.g. "100"), not as a
# number, matching the original workbook's shared-string type. Force text
# formatting before writing, then restore the default "Normal" style so no
# stray number-format styling is left behind on the cells.
$ws.Range("B2:B10").NumberFormat = "@"

$ws.Range("A2").Value = "2244-BLSM-L"
$ws.Range("B2").Value = "100|200"

$ws.Range("A3").Value = "2244-BLSM-M"
$ws.Range("B3").Value = "100|200"

$ws.Range("A4").Value = "SLKDRM-CLK-52-L"
$ws.Range("B4").Value = "100"

$ws.Range("A5").Value = "SLKDRM-CLK-52-S"
$ws.Range("B5").Value = "100"

$ws.Range("A6").Value = "2244-BLSM-S"
$ws.Range("B6").Value = "100|200"

$ws.Range("A7").Value = "SLKDRM-CLK-52-XL"
$ws.Range("B7").Value = "100"

$ws.Range("A8").Value = "SLKDRM-CLK-52-M"
$ws.Range("B8").Value = "100"

$ws.Range("A9").Value = "2244-BLSM-XL"
$ws.Range("B9").Value = "100|200"

$ws.Range("A10").Value = "SLKDRM-CLK-03-XL"
$ws.Range("B10").Value = "100"

$ws.Range("B2:B10").Style = "Normal"
